# Natmi following Dr Hou advice
#
# Adds a new "sCs" cluster into the LR-pairs table (Col9a2 -> Mag) and
# re-derives the specificity table for all Sending/Target cluster
# combinations between the existing "FAPs" cluster and the new "sCs"
# cluster, expanding the table from a single data row to four data rows
# (rows 2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Col9a2/Mag -> M2
$ws.Cells.Item(2, 1).Value  = "FAPs"
$ws.Cells.Item(2, 2).Value  = "Col9a2"
$ws.Cells.Item(2, 3).Value  = "Mag"
$ws.Cells.Item(2, 4).Value  = "M2"
$ws.Cells.Item(2, 5).Value  = 3
$ws.Cells.Item(2, 6).Value  = 1
$ws.Cells.Item(2, 7).Value  = 0.5218163333333333
$ws.Cells.Item(2, 8).Value  = 1.565449
$ws.Cells.Item(2, 9).Value  = 0.9834142768387242
$ws.Cells.Item(2, 10).Value = 0.9834142768387242
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.7885686666666668
$ws.Cells.Item(2, 14).Value = 2.365706
$ws.Cells.Item(2, 15).Value = 0.4566863346753138
$ws.Cells.Item(2, 16).Value = 0.4566863346753137
$ws.Cells.Item(2, 17).Value = 0.4114880102215556
$ws.Cells.Item(2, 18).Value = 3.703392091994001
$ws.Cells.Item(2, 19).Value = 0.4491118615568513
$ws.Cells.Item(2, 20).Value = 0.4491118615568512

# Row 3: FAPs -> Col9a2/Mag -> sCs
$ws.Cells.Item(3, 1).Value  = "FAPs"
$ws.Cells.Item(3, 2).Value  = "Col9a2"
$ws.Cells.Item(3, 3).Value  = "Mag"
$ws.Cells.Item(3, 4).Value  = "sCs"
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 6).Value  = 1
$ws.Cells.Item(3, 7).Value  = 0.5218163333333333
$ws.Cells.Item(3, 8).Value  = 1.565449
$ws.Cells.Item(3, 9).Value  = 0.9834142768387242
$ws.Cells.Item(3, 10).Value = 0.9834142768387242
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.9381496666666668
$ws.Cells.Item(3, 14).Value = 2.814449
$ws.Cells.Item(3, 15).Value = 0.5433136653246862
$ws.Cells.Item(3, 16).Value = 0.5433136653246862
$ws.Cells.Item(3, 17).Value = 0.4895418191778889
$ws.Cells.Item(3, 18).Value = 4.405876372601001
$ws.Cells.Item(3, 19).Value = 0.5343024152818728
$ws.Cells.Item(3, 20).Value = 0.5343024152818728

# Row 4: sCs -> Col9a2/Mag -> M2
$ws.Cells.Item(4, 1).Value  = "sCs"
$ws.Cells.Item(4, 2).Value  = "Col9a2"
$ws.Cells.Item(4, 3).Value  = "Mag"
$ws.Cells.Item(4, 4).Value  = "M2"
$ws.Cells.Item(4, 5).Value  = 1
$ws.Cells.Item(4, 6).Value  = 0.3333333333333333
$ws.Cells.Item(4, 7).Value  = 0.008800666666666667
$ws.Cells.Item(4, 8).Value  = 0.026402
$ws.Cells.Item(4, 9).Value  = 0.01658572316127577
$ws.Cells.Item(4, 10).Value = 0.01658572316127577
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7885686666666668
$ws.Cells.Item(4, 14).Value = 2.365706
$ws.Cells.Item(4, 15).Value = 0.4566863346753138
$ws.Cells.Item(4, 16).Value = 0.4566863346753137
$ws.Cells.Item(4, 17).Value = 0.006939929979111113
$ws.Cells.Item(4, 18).Value = 0.062459369812
$ws.Cells.Item(4, 19).Value = 0.007574473118462492
$ws.Cells.Item(4, 20).Value = 0.007574473118462489

# Row 5: sCs -> Col9a2/Mag -> sCs
$ws.Cells.Item(5, 1).Value  = "sCs"
$ws.Cells.Item(5, 2).Value  = "Col9a2"
$ws.Cells.Item(5, 3).Value  = "Mag"
$ws.Cells.Item(5, 4).Value  = "sCs"
$ws.Cells.Item(5, 5).Value  = 1
$ws.Cells.Item(5, 6).Value  = 0.3333333333333333
$ws.Cells.Item(5, 7).Value  = 0.008800666666666667
$ws.Cells.Item(5, 8).Value  = 0.026402
$ws.Cells.Item(5, 9).Value  = 0.01658572316127577
$ws.Cells.Item(5, 10).Value = 0.01658572316127577
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.9381496666666668
$ws.Cells.Item(5, 14).Value = 2.814449
$ws.Cells.Item(5, 15).Value = 0.5433136653246862
$ws.Cells.Item(5, 16).Value = 0.5433136653246862
$ws.Cells.Item(5, 17).Value = 0.008256342499777778
$ws.Cells.Item(5, 18).Value = 0.074307082498
$ws.Cells.Item(5, 19).Value = 0.009011250042813283
$ws.Cells.Item(5, 20).Value = 0.009011250042813279
